$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.603556632995605
$ws.Range("B1").Value = 6.665385723114014
$ws.Range("C1").Value = 6.09206485748291
$ws.Range("D1").Value = 4.864583015441895
$ws.Range("E1").Value = 2.166341304779053
